# Fruta / hortaliza, semanal
# Insert a new week's worth of observations (2 rows) for
# "Vega Monumental Concepción - Kiwi" ahead of the existing records,
# pushing the previous data down by two rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at the top of the data block (row 165),
# shifting the former rows 165:184 down to 167:186.
$ws.Rows("165:166").Insert()

# --- Row 165: Kiwi Hayward "Primera" for the new week ---
$ws.Range("A165").Value = 11
$ws.Range("B165").Value = "Vega Monumental Concepción"
$ws.Range("C165").Value = "Bíobío"
$ws.Range("D165").Value = 44769
$ws.Range("E165").Value = 8
$ws.Range("F165").Value = "Fruta"
$ws.Range("G165").Value = 100101
$ws.Range("H165").Value = "Berries"
$ws.Range("I165").Value = 100101007
$ws.Range("J165").Value = "Kiwi"
$ws.Range("K165").Value = "Hayward"
$ws.Range("L165").Value = "Primera"
$ws.Range("M165").Value = 200
$ws.Range("N165").Value = 6000
$ws.Range("O165").Value = 7000
$ws.Range("P165").Value = 6500
$ws.Range("Q165").Value = "$/bandeja 18 kilos"
$ws.Range("R165").Value = "Región de O'Higgins"
$ws.Range("S165").Value = 361
$ws.Range("T165").Value = 18

# --- Row 166: Kiwi Hayward "Segunda" for the new week ---
$ws.Range("A166").Value = 11
$ws.Range("B166").Value = "Vega Monumental Concepción"
$ws.Range("C166").Value = "Bíobío"
$ws.Range("D166").Value = 44769
$ws.Range("E166").Value = 8
$ws.Range("F166").Value = "Fruta"
$ws.Range("G166").Value = 100101
$ws.Range("H166").Value = "Berries"
$ws.Range("I166").Value = 100101007
$ws.Range("J166").Value = "Kiwi"
$ws.Range("K166").Value = "Hayward"
$ws.Range("L166").Value = "Segunda"
$ws.Range("M166").Value = 200
$ws.Range("N166").Value = 5000
$ws.Range("O166").Value = 5500
$ws.Range("P166").Value = 5250
$ws.Range("Q166").Value = "$/bandeja 18 kilos"
$ws.Range("R166").Value = "Región de O'Higgins"
$ws.Range("S166").Value = 292
$ws.Range("T166").Value = 18
